$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C got (re)sized - closest achievable value to the authored 9.125 ---
$ws.Columns.Item(3).ColumnWidth = 8.3

# --- New row 7 tail cells (X7 / Y7) ---
$ws.Range("X7").Value = 0.42000000000000171
$ws.Range("Y7").Value = "Up"

# --- New row 8 ---
$ws.Range("A8").Value = 42649.879930555559
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = "Strong Buy"
$ws.Range("D8:O8").Value = 0
$ws.Range("P8").Value = "Random"
$ws.Range("Q8").Value = 47.321424984051369
$ws.Range("R8").Value = 0
$ws.Range("S8").NumberFormat = "0.00%"
$ws.Range("S8").Value = 0.0834
$ws.Range("T8").NumberFormat = "0.00%"
$ws.Range("T8").Value = -0.0062
$ws.Range("U8").Value = 2.31
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0
